# Add "2022-Q1" fund-detail sheet (inserted right after "2021-Q4", before "总计"),
# and add a matching "2022-Q1" summary row at the top of the "总计" sheet's data.
#
# NOTE: worksheet object handles obtained by index/name become stale once the
# sheet collection is restructured (e.g. via Worksheets.Add); any sheet whose
# position shifts as a result must be re-fetched afterwards. So we finish all
# of the sheet-collection surgery first, and only then (re-)resolve "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet right after "2021-Q4" (i.e. right before "总计")
#    and name it "2022-Q1". Copy the cell formatting (borders/bold header,
#    index-column style, …) from the "2021-Q4" sheet so it matches the look
#    of the other per-quarter fund sheets.
# ---------------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Item("2021-Q4")

$newSheet = $wb.Worksheets.Add($null, $q4Sheet)
$newSheet.Name = "2022-Q1"

# Columns B-G hold text (fund code, name, scale, position %, …) even when the
# text looks numeric, so force text formatting up front; do this BEFORE the
# border/bold formatting copy below so the header/index-column look (copied
# from "2021-Q4") still wins on the cells it actually covers.
$newSheet.Range("B2:G5").NumberFormat = "@"

# Header row look (bold + border), copied from the "2021-Q4" header.
$q4Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

# Index column (A) look (bold + border), copied from "2021-Q4"; its data
# block is only 2 rows long but the new sheet needs 4, so paste row-by-row.
$q4Sheet.Range("A2").Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122)   # xlPasteFormats

# Match page margins used by the other non-first sheets.
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------------
# 2. Fill in header row.
# ---------------------------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------------
# 3. Fill in the fund holding rows (A = index, H = rank stay numeric; B-G
#    were already forced to text format above).
# ---------------------------------------------------------------------------
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "501305"
$newSheet.Range("C2").Value = "汇添富中证港股通高股息投资指数（LOF）A"
$newSheet.Range("D2").Value = "1.59"
$newSheet.Range("E2").Value = "93.08"
$newSheet.Range("F2").Value = "3.52"
$newSheet.Range("G2").Value = "0.0560"
$newSheet.Range("H2").Value = 6

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "501306"
$newSheet.Range("C3").Value = "汇添富中证港股通高股息投资指数（LOF）C"
$newSheet.Range("D3").Value = "0.21"
$newSheet.Range("E3").Value = "93.08"
$newSheet.Range("F3").Value = "3.52"
$newSheet.Range("G3").Value = "0.0074"
$newSheet.Range("H3").Value = 6

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "501307"
$newSheet.Range("C4").Value = "银河中证沪港深高股息指数（LOF）A"
$newSheet.Range("D4").Value = "0.19"
$newSheet.Range("E4").Value = "91.35"
$newSheet.Range("F4").Value = "1.52"
$newSheet.Range("G4").Value = "0.0029"
$newSheet.Range("H4").Value = 7

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "501308"
$newSheet.Range("C5").Value = "银河中证沪港深高股息指数（LOF）C"
$newSheet.Range("D5").Value = "0.01"
$newSheet.Range("E5").Value = "91.35"
$newSheet.Range("F5").Value = "1.52"
$newSheet.Range("G5").Value = "0.0002"
$newSheet.Range("H5").Value = 7

# ---------------------------------------------------------------------------
# 4. Prepend a "2022-Q1" row to the "总计" sheet (new row 2), pushing the
#    existing rows down, and renumber the index column (A) for every row.
#    Re-resolve the "总计" worksheet now that the sheet collection has
#    settled, since its index shifted (4 -> 5) when we inserted "2022-Q1".
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows.Item(2).Insert()

# The inserted row inherited formatting from the row below it; re-apply the
# correct per-column look (bold/border index-col style on A, plain on B:D)
# by pulling formats from row 3, which still carries the original styling.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)    # xlPasteFormats
$totalSheet.Range("B3:D3").Copy()
$totalSheet.Range("B2:D2").PasteSpecial(-4122) # xlPasteFormats

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.07000000000000001

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3

# ---------------------------------------------------------------------------
# 5. Restore the originally active sheet/selection (inserting + renaming a
#    sheet shifts Excel's active-tab focus onto it as a side effect).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q1").Activate()
